# Fix situation_cloture issue <contrat resilie> and update checkContratsAv function
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The contract for "Karami abdelilah" (old row 3, CIN BB779645) was resiliated
# and must be removed entirely; the rows below it shift up by one.
$ws.Rows(3).Delete()

# Row 2 now reflects a different, corrected contract/tenant record.
$ws.Range("A2").Value = "794/DR KESH"
$ws.Range("C2").Value = "BG1949"
$ws.Range("D2").Value = "Ahmed Test"
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = 8500

# The trailing totals row (now row 7 after the row shift) must be recomputed
# to reflect the removed contract.
$ws.Range("I7").Value = 10000
$ws.Range("K7").Value = 1500
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = 26800
